# Changes of 23rd March 2022
# Update the ShipmentTracking (P), ActualRate (Q) and Result (R) columns for the
# shipment rows (2-25) on Sheet1.
#
# Tracking numbers and dollar amounts look numeric, so a leading apostrophe is used
# when assigning them (exactly like typing '320018110083 into Excel) to force them to
# stay text, matching the original t="s" shared-string cells, instead of silently
# becoming numeric values. The Style reset right after clears the "quote prefix" flag
# that the apostrophe entry leaves behind, so no stray style/number-format is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @{ Row = 2; P = '320018110083'; Q = ''; R = '' },
    @{ Row = 3; P = '320018110094'; Q = ''; R = '' },
    @{ Row = 4; P = '320018110120'; Q = ''; R = '' },
    @{ Row = 5; P = '320018110142'; Q = ''; R = '' },
    @{ Row = 6; P = '320018110186'; Q = ''; R = '' },
    @{ Row = 7; P = '320018110201'; Q = ''; R = '' },
    @{ Row = 8; P = '320018110234'; Q = ''; R = '' },
    @{ Row = 9; P = '320018110256'; Q = ''; R = '' },
    @{ Row = 10; P = '320018110289'; Q = ''; R = '' },
    @{ Row = 11; P = '320018110304'; Q = ''; R = '' },
    @{ Row = 12; P = '320018110348'; Q = ''; R = '' },
    @{ Row = 13; P = '320018110360'; Q = ''; R = '' },
    @{ Row = 14; P = '320018110392'; Q = ''; R = '' },
    @{ Row = 15; P = '320018110418'; Q = ''; R = '' },
    @{ Row = 16; P = '320018110440'; Q = ''; R = '' },
    @{ Row = 17; P = '320018110462'; Q = ''; R = '' },
    @{ Row = 18; P = '320018110500'; Q = '$53.93'; R = 'FAIL' },
    @{ Row = 19; P = '320018110521'; Q = ''; R = '' },
    @{ Row = 20; P = '320018110554'; Q = '$85.66'; R = 'FAIL' },
    @{ Row = 21; P = '320018110576'; Q = ''; R = '' },
    @{ Row = 22; P = '320018110602'; Q = '$195.48'; R = '' },
    @{ Row = 23; P = '320018110613'; Q = '$439.28'; R = 'PASS' },
    @{ Row = 24; P = '320018110624'; Q = ''; R = '' },
    @{ Row = 25; P = '320018110635'; Q = ''; R = '' }
)

foreach ($item in $rows) {
    $pCell = $ws.Range('P' + $item.Row)
    $pCell.Value = "'" + $item.P
    $pCell.Style = 'Normal'

    if ($item.Q -ne '') {
        $qCell = $ws.Range('Q' + $item.Row)
        $qCell.Value = "'" + $item.Q
        $qCell.Style = 'Normal'
    }

    if ($item.R -ne '') {
        $rCell = $ws.Range('R' + $item.Row)
        $rCell.Value = $item.R
    }
}

Write-Output 'Updated ShipmentTracking/ActualRate/Result for rows 2-25 on Sheet1'
